# Auto-generated edit script: updates market-price-derived columns (H-N)
# for specific Leve rows across multiple sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 67849.53
$ws.Range("I70").Value = 817.5
$ws.Range("K70").Value = 2452.5
$ws.Range("M70").Value = -2182.5
# Row 73
$ws.Range("H73").Value = 67849.53
$ws.Range("I73").Value = 817.5
$ws.Range("K73").Value = 2452.5
$ws.Range("M73").Value = -1516.5
# Row 98
$ws.Range("H98").Value = 2024.5
$ws.Range("J98").Value = 4701
$ws.Range("L98").Value = 4701
$ws.Range("N98").Value = -7697
# Row 107
$ws.Range("H107").Value = 50554.094
$ws.Range("I107").Value = 62037.824
$ws.Range("K107").Value = 62037.824
$ws.Range("M107").Value = -60117.824
# Row 108
$ws.Range("H108").Value = 68400
$ws.Range("J108").Value = 68400
$ws.Range("L108").Value = 68400
$ws.Range("N108").Value = -76080
# Row 122
$ws.Range("H122").Value = 2024.5
$ws.Range("J122").Value = 4701
$ws.Range("L122").Value = 14103
$ws.Range("N122").Value = -19003
# Row 132
$ws.Range("H132").Value = 5005.1816
$ws.Range("I132").Value = 5431.25
$ws.Range("J132").Value = 3869
$ws.Range("K132").Value = 16293.75
$ws.Range("L132").Value = 11607
$ws.Range("M132").Value = -13763.75
$ws.Range("N132").Value = -16667
# Row 137
$ws.Range("H137").Value = 3229.4075
$ws.Range("I137").Value = 1976.7333
$ws.Range("J137").Value = 4795.25
$ws.Range("K137").Value = 5930.199900000001
$ws.Range("L137").Value = 14385.75
$ws.Range("M137").Value = -3380.199900000001
$ws.Range("N137").Value = -19485.75
# Row 138
$ws.Range("H138").Value = 5503.636
$ws.Range("J138").Value = 6715.0415
$ws.Range("L138").Value = 20145.1245
$ws.Range("N138").Value = -30425.1245

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 5105.88
$ws.Range("I45").Value = 1724.1428
$ws.Range("K45").Value = 1724.1428
$ws.Range("M45").Value = -1347.1428
# Row 61
$ws.Range("H61").Value = 3029.8484
$ws.Range("I61").Value = 1662.7273
$ws.Range("J61").Value = 5764.091
$ws.Range("K61").Value = 1662.7273
$ws.Range("L61").Value = 5764.091
$ws.Range("M61").Value = -1450.7273
$ws.Range("N61").Value = -6188.091
# Row 63
$ws.Range("H63").Value = 8476.25
$ws.Range("I63").Value = 3905
$ws.Range("K63").Value = 3905
$ws.Range("M63").Value = -3219
# Row 64
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
# Row 66
$ws.Range("H66").Value = 8476.25
$ws.Range("I66").Value = 3905
$ws.Range("K66").Value = 19525
$ws.Range("M66").Value = -16093
# Row 67
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
# Row 80
$ws.Range("H80").Value = 63950
$ws.Range("J80").Value = 63950
$ws.Range("L80").Value = 63950
$ws.Range("N80").Value = -65946
# Row 82
$ws.Range("H82").Value = 39999
$ws.Range("J82").Value = 39999
$ws.Range("L82").Value = 39999
$ws.Range("N82").Value = -40721
# Row 83
$ws.Range("H83").Value = 63950
$ws.Range("J83").Value = 63950
$ws.Range("L83").Value = 191850
$ws.Range("N83").Value = -201834
# Row 85
$ws.Range("H85").Value = 39999
$ws.Range("J85").Value = 39999
$ws.Range("L85").Value = 39999
$ws.Range("N85").Value = -42495
# Row 86
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32372
# Row 88
$ws.Range("H88").Value = 1850.7778
$ws.Range("I88").Value = 865.3333
$ws.Range("J88").Value = 2343.5
$ws.Range("K88").Value = 865.3333
$ws.Range("L88").Value = 2343.5
$ws.Range("M88").Value = -459.3333
$ws.Range("N88").Value = -3155.5
# Row 89
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 90000
$ws.Range("N89").Value = -101856
# Row 91
$ws.Range("H91").Value = 1850.7778
$ws.Range("I91").Value = 865.3333
$ws.Range("J91").Value = 2343.5
$ws.Range("K91").Value = 865.3333
$ws.Range("L91").Value = 2343.5
$ws.Range("M91").Value = 538.6667
$ws.Range("N91").Value = -5151.5
# Row 122
$ws.Range("H122").Value = 4800.5
$ws.Range("I122").Value = 3252
$ws.Range("J122").Value = 5832.8335
$ws.Range("K122").Value = 9756
$ws.Range("L122").Value = 17498.5005
$ws.Range("M122").Value = -7306
$ws.Range("N122").Value = -22398.5005
# Row 136
$ws.Range("H136").Value = 3029.8484
$ws.Range("I136").Value = 1662.7273
$ws.Range("J136").Value = 5764.091
$ws.Range("K136").Value = 4988.1819
$ws.Range("L136").Value = 17292.273
$ws.Range("M136").Value = -2438.1819
$ws.Range("N136").Value = -22392.273

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4590.1816
$ws.Range("I105").Value = 2899.6
$ws.Range("J105").Value = 5999
$ws.Range("K105").Value = 2899.6
$ws.Range("L105").Value = 5999
$ws.Range("M105").Value = -1152.6
$ws.Range("N105").Value = -9493

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3137.4546
$ws.Range("I31").Value = 1010.0909
$ws.Range("K31").Value = 1010.0909
$ws.Range("M31").Value = -715.0909
# Row 34
$ws.Range("H34").Value = 3137.4546
$ws.Range("I34").Value = 1010.0909
$ws.Range("K34").Value = 1010.0909
$ws.Range("M34").Value = -808.0909
# Row 62
$ws.Range("H62").Value = 3962
$ws.Range("J62").Value = 3943.5
$ws.Range("L62").Value = 3943.5
$ws.Range("N62").Value = -5191.5
# Row 65
$ws.Range("H65").Value = 3962
$ws.Range("J65").Value = 3943.5
$ws.Range("L65").Value = 19717.5
$ws.Range("N65").Value = -25957.5
# Row 99
$ws.Range("H99").Value = 4470.3335
$ws.Range("I99").Value = 3461.8572
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 3461.8572
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -1963.8572
$ws.Range("N99").Value = -10996
# Row 105
$ws.Range("H105").Value = 1895.5
$ws.Range("I105").Value = 1772.7778
$ws.Range("K105").Value = 1772.7778
$ws.Range("M105").Value = -25.77780000000007
# Row 122
$ws.Range("H122").Value = 3811.5
$ws.Range("I122").Value = 1264.25
$ws.Range("J122").Value = 4830.4
$ws.Range("K122").Value = 3792.75
$ws.Range("L122").Value = 14491.2
$ws.Range("M122").Value = -1342.75
$ws.Range("N122").Value = -19391.2
# Row 126
$ws.Range("H126").Value = 4470.3335
$ws.Range("I126").Value = 3461.8572
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 10385.5716
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -7915.571599999999
$ws.Range("N126").Value = -28940
# Row 132
$ws.Range("H132").Value = 2946.6155
$ws.Range("I132").Value = 2726.15
$ws.Range("K132").Value = 8178.450000000001
$ws.Range("M132").Value = -5648.450000000001
# Row 134
$ws.Range("H134").Value = 458559.78
$ws.Range("I134").Value = 3878.2307
$ws.Range("K134").Value = 11634.6921
$ws.Range("M134").Value = -9099.6921

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 4.5
$ws.Range("J46").Value = 8
$ws.Range("L46").Value = 24
$ws.Range("N46").Value = -206

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2004777.2
$ws.Range("I80").Value = 3336970.2
$ws.Range("J80").Value = 1433837.4
$ws.Range("K80").Value = 3336970.2
$ws.Range("L80").Value = 1433837.4
$ws.Range("M80").Value = -3335972.2
$ws.Range("N80").Value = -1435833.4
# Row 83
$ws.Range("H83").Value = 2004777.2
$ws.Range("I83").Value = 3336970.2
$ws.Range("J83").Value = 1433837.4
$ws.Range("K83").Value = 16684851
$ws.Range("L83").Value = 7169187
$ws.Range("M83").Value = -16679859
$ws.Range("N83").Value = -7179171

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 362800
$ws.Range("I7").Value = 631001.25
$ws.Range("J7").Value = 5198.3335
$ws.Range("K7").Value = 631001.25
$ws.Range("L7").Value = 5198.3335
$ws.Range("M7").Value = -630889.25
$ws.Range("N7").Value = -5422.3335
# Row 22
$ws.Range("H22").Value = 703.8261
$ws.Range("I22").Value = 762.85
$ws.Range("J22").Value = 310.33334
$ws.Range("K22").Value = 762.85
$ws.Range("L22").Value = 310.33334
$ws.Range("M22").Value = -467.85
$ws.Range("N22").Value = -900.33334
# Row 27
$ws.Range("H27").Value = 703.8261
$ws.Range("I27").Value = 762.85
$ws.Range("J27").Value = 310.33334
$ws.Range("K27").Value = 762.85
$ws.Range("L27").Value = 310.33334
$ws.Range("M27").Value = -655.85
$ws.Range("N27").Value = -524.33334
# Row 38
$ws.Range("H38").Value = 29999
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
# Row 40
$ws.Range("H40").Value = 1117078.2
$ws.Range("I40").Value = 2004140.8
$ws.Range("J40").Value = 8250
$ws.Range("K40").Value = 2004140.8
$ws.Range("L40").Value = 8250
$ws.Range("M40").Value = -2004004.8
$ws.Range("N40").Value = -8522
# Row 126
$ws.Range("H126").Value = 362800
$ws.Range("I126").Value = 631001.25
$ws.Range("J126").Value = 5198.3335
$ws.Range("K126").Value = 1893003.75
$ws.Range("L126").Value = 15595.0005
$ws.Range("M126").Value = -1890533.75
$ws.Range("N126").Value = -20535.0005
# Row 132
$ws.Range("H132").Value = 8749.5
$ws.Range("J132").Value = 9090.362999999999
$ws.Range("L132").Value = 27271.089
$ws.Range("N132").Value = -32331.089

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 10007
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 10007
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -10291
# Row 126
$ws.Range("H126").Value = 3555.2632
$ws.Range("I126").Value = 3111.5557
$ws.Range("K126").Value = 9334.667099999999
$ws.Range("M126").Value = -6864.667099999999
